# Update "想去人数" (number of people interested) values on the
# "展览" (Exhibition) and "全部类型" (All types) sheets to reflect
# the latest scrape results.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 7449
$wsExhibition.Range("F5").Value = 445
$wsExhibition.Range("F6").Value = 4008
$wsExhibition.Range("F7").Value = 320
$wsExhibition.Range("F8").Value = 564
$wsExhibition.Range("F10").Value = 642
$wsExhibition.Range("F11").Value = 124

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7449
$wsAll.Range("F7").Value = 445
$wsAll.Range("F8").Value = 4008
$wsAll.Range("F9").Value = 320
$wsAll.Range("F10").Value = 564
$wsAll.Range("F12").Value = 642
$wsAll.Range("F14").Value = 124
